# "fix A4 A6; Add B1 B2"
#
# 1) Fix the "Uptick_7000-2" typo (wrong capitalisation) that is shown on
#    sheets A4 and A6 (cell D2 on both) -> "uptick_7000-2".
# 2) Populate the B1 and B2 sheets (previously placeholder copies of the
#    B5/B6 "TxHash" sheets) with their own TxHash values.
# 3) Make B2 the active sheet/tab (it was A20 before the edit).

$wb = $excel.ActiveWorkbook

# --- 1) Fix A4 / A6 -------------------------------------------------------
$wsA4 = $wb.Worksheets.Item("A4")
$wsA4.Range("D2").Value = "uptick_7000-2"

$wsA6 = $wb.Worksheets.Item("A6")
$wsA6.Range("D2").Value = "uptick_7000-2"

# --- 2) Add B1 / B2 content ------------------------------------------------
$wsB1 = $wb.Worksheets.Item("B1")
$wsB1.Range("A2").Value = "F705B0C56E3ED3F5AB6F948C0E682EA2B8567D3BE3A46082C743B616E74AD286"
$wsB1.Range("A3").Value = "E7B6B928963AA3FD27522852DE0F7812D211837DD7CEA9D3DFD7A212E150D678"

$wsB2 = $wb.Worksheets.Item("B2")
$wsB2.Range("A2").Value = "B126D40FA0DE01443BD35F9B417BFFC8AE817BAC16315650063A56CC4F931A88"
$wsB2.Range("A3").Value = "721E3218B7D72A06E21BFC9210DA9A750B29DCF7D27A657ABBAB962F34593225"

# --- 3) Make B2 the active sheet ------------------------------------------
$wsB2.Activate()
